# Add two new columns to the worksheet: I ("I0") and J ("IF").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy formatting (font/border/alignment) from the
# existing header cell H1 so the new headers match the table's look.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-32: (I value, J value) pairs, in row order.
$data = @(
    @(10, 10),
    @(4, 5),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(1, 1),
    @(8, 8),
    @(5, 5),
    @(7, 7),
    @(8, 9),
    @(4, 5),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(11, 11),
    @(9, 9),
    @(2, 3),
    @(9, 9),
    @(11, 11),
    @(8, 8),
    @(7, 7),
    @(5, 5),
    @(5, 5),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $pair = $data[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
